# Update activity stats (runs, balls, fours, sixes) for Shivam Dube rows.
# The underlying data rows (2-10) have been reshuffled/updated to reflect
# the latest match activity, per the "updated activity till excel form" commit.
# Only the cells whose value actually changes are touched (row 6 - and a
# handful of other cells whose new value equals the old one - are left
# completely alone so no incidental formatting churn is introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> column index for runs/balls/fours/sixes.
$cols = @{ "C" = 3; "D" = 4; "E" = 5; "F" = 6 }

# Only the cells whose value actually changes, row -> column -> new value.
$changes = @{
    2  = @{ "C" = "8";  "D" = "13" }
    3  = @{ "C" = "17"; "D" = "11"; "E" = "2"; "F" = "1" }
    4  = @{ "C" = "2";  "D" = "6";  "E" = "0"; "F" = "0" }
    5  = @{ "C" = "12"; "D" = "12"; "E" = "0" }
    7  = @{ "C" = "22"; "D" = "14"; "E" = "2"; "F" = "1" }
    8  = @{ "C" = "27"; "D" = "10"; "E" = "1"; "F" = "3" }
    9  = @{ "C" = "7";  "D" = "8";  "F" = "0" }
    10 = @{ "C" = "11"; "D" = "12"; "F" = "1" }
}

foreach ($row in $changes.Keys) {
    $rowChanges = $changes[$row]
    foreach ($col in $rowChanges.Keys) {
        $cell = $ws.Cells.Item($row, $cols[$col])

        # Preserve the text (string) storage of these numeric-looking values,
        # matching the source file's original number-stored-as-text convention.
        $cell.NumberFormat = "@"
        $cell.Value = $rowChanges[$col]
    }
}
